# Update "想去人数" (interest count) figures on the "展览" and "全部类型"
# sheets to reflect newly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1149
$wsExpo.Range("F4").Value = 2605
$wsExpo.Range("F5").Value = 227

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1149
$wsAll.Range("F6").Value = 2605
$wsAll.Range("F8").Value = 227
